$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.005.10'
$ws.Range("E2").Value = '  -0.49%  '

$ws.Range("D3").Value = '1.859.22'
$ws.Range("E3").Value = '  -1.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.24'
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5092'
$ws.Range("E7").Value = '  +0.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3838'
$ws.Range("E8").Value = '  -0.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08223'
$ws.Range("E9").Value = '  -9.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.112'
$ws.Range("E10").Value = '  -1.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.53'
$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.213'
$ws.Range("E12").Value = '  -2.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.58'
$ws.Range("E13").Value = '  -1.29%  '

$ws.Range("D14").Value = '1.859.92'
$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.257'
$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("E17").Value = '  -1.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.79'
$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06651'
$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.69'
$ws.Range("E20").Value = '  -3.27%  '

$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.017'
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").Value = '28.036.02'
$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("E24").Value = '  -3.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.236'
$ws.Range("E25").Value = '  -1.31%  '

$ws.Range("D26").Value = '2.073.18'
$ws.Range("E26").Value = '  -0.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.517'
$ws.Range("E27").Value = '  -1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.20'
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("E29").Value = '  -1.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.88'
$ws.Range("E30").Value = '  -1.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1059'
$ws.Range("E31").Value = '  -0.40%  '

$ws.Range("E32").Value = '  -2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.927'
$ws.Range("E33").Value = '  +5.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.595'
$ws.Range("E34").Value = '  -0.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.387'
$ws.Range("E35").Value = '  -1.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06517'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02413'
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2173'
$ws.Range("E38").Value = '  -1.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6553'
$ws.Range("E39").Value = '  +1.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.197'
$ws.Range("E40").Value = '  -1.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.003'
$ws.Range("E41").Value = '  +1.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.222'
$ws.Range("E42").Value = '  -5.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.17'
$ws.Range("E43").Value = '  -3.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6132'
$ws.Range("E44").Value = '  +1.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.12'
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.656'
$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("E48").Value = '  +0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.207'
$ws.Range("E49").Value = '  -2.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.17'
$ws.Range("E50").Value = '  -1.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.36'
$ws.Range("E51").Value = '  -1.64%  '
